# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets, per the upstream data refresh
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row => new value, for the "展览" sheet (sheet1)
$exhibitUpdates = @{
    2  = 35
    4  = 94
    6  = 560
    7  = 1754
    8  = 43
    10 = 141
    11 = 1934
    13 = 338
    14 = 447
    16 = 281
    17 = 205
    19 = 20
    21 = 42
    22 = 55
    23 = 1046
    25 = 321
    26 = 174
    27 = 258
    28 = 293
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row => new value, for the "全部类型" sheet (sheet4)
$allUpdates = @{
    2  = 35
    4  = 94
    6  = 560
    7  = 1754
    9  = 43
    11 = 141
    12 = 1934
    14 = 338
    15 = 447
    17 = 281
    18 = 205
    20 = 20
    22 = 42
    23 = 55
    24 = 1046
    26 = 321
    27 = 174
    28 = 258
    29 = 293
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
